# fire truck report, JWT on
# Rework the "passenger car" operational-log header row: several header
# labels are clarified with units (", л" / ", км"), the fuel-refill /
# return-fuel labels are renamed, the driver-name label is tidied up, and
# two new "distance by city / by region" columns (F, G) are introduced
# (disambiguating them from the existing "Город"/"Область" sub-headers
# under "Израсходовано ГСМ").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 header text updates -------------------------------------------
# Order matters: the workbook's shared-string table is appended to in the
# exact sequence cells are (re)written, so we follow the same edit order
# the original author used.
$ws.Range("C4").Value = "Наличие ГСМ перед выездом, л"
$ws.Range("D4").Value = "Показания спидометра перед выездом, км"
$ws.Range("L4").Value = "Получено ГСМ, л"
$ws.Range("M4").Value = "Наличие ГСМ при возвращении, л"
$ws.Range("N4").Value = "Показания спидометра при возвращении, км"
$ws.Range("O4").Value = "Экономия, л"
$ws.Range("P4").Value = "Перерасход, л"
$ws.Range("B4").Value = "Фамилия И. О. водителя"
$ws.Range("F4").Value = "Пройдено км по городу"
$ws.Range("G4").Value = "Пройдено км по области"

# --- Column width tweaks --------------------------------------------------
# New columns F/G need explicit widths now that they carry real headers;
# column M widens slightly to fit its longer label.
$ws.Columns.Item(6).ColumnWidth = 10.585
$ws.Columns.Item(7).ColumnWidth = 10.085
$ws.Columns.Item(13).ColumnWidth = 14.42

# --- Title row font tidy-up ------------------------------------------------
# The banner row (A1:P1) was set in a stray 10pt weight; bump it to the
# same 11pt bold Times New Roman used by every other header on the sheet.
$ws.Range("A1:P1").Font.Size = 11

# --- Selection ------------------------------------------------------------
# Land the selection/scroll position on the subtitle row instead of the
# previously-saved far-right cell.
$ws.Range("A3:P3").Select()
